$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.481909168789134
$ws.Range("B2").Value = -4.230807130045484

$ws.Range("A3").Value = -0.4649429187685051
$ws.Range("B3").Value = 0.7649353906709653

$ws.Range("B4").Value = -3.18222582427557

$ws.Range("A5").Value = 0.7242133972617233
$ws.Range("B5").Value = 0.5215895993778858

$ws.Range("A6").Value = -0.8192272649316585
$ws.Range("B6").Value = -1.970377164371514

$ws.Range("A7").Value = -0.09385660014301761
$ws.Range("B7").Value = -0.9361813794003838

$ws.Range("B8").Value = 0.5227215113784796

$ws.Range("A9").Value = 0.3310252182837974
$ws.Range("B9").Value = 0.7865897687033396

$ws.Range("A10").Value = -0.1876020222794525
$ws.Range("B10").Value = -2.036259574405221

$ws.Range("A11").Value = 0.2927036682248865
$ws.Range("B11").Value = -0.4420696700903583
